$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 151, shifting existing rows 151-157 down to 152-158
$ws.Rows.Item(151).Insert()

# Populate the new row 151 with the weekly price entry (same series as row 150,
# new week's date + same price bracket)
$ws.Cells.Item(151, 1).Value = 11
$ws.Cells.Item(151, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(151, 3).Value = "Bíobío"
$ws.Cells.Item(151, 4).Value = 44826
$ws.Cells.Item(151, 5).Value = 8
$ws.Cells.Item(151, 6).Value = 100112043
$ws.Cells.Item(151, 7).Value = "Pepino ensalada"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 100
$ws.Cells.Item(151, 11).Value = 17000
$ws.Cells.Item(151, 12).Value = 18000
$ws.Cells.Item(151, 13).Value = 17500
$ws.Cells.Item(151, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(151, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(151, 16).Value = 292
$ws.Cells.Item(151, 17).Value = 60
$ws.Cells.Item(151, 18).Value = "Hortaliza"

# Match the date cell style used by the other rows in column D (row 150 as reference)
$ws.Cells.Item(151, 4).NumberFormat = $ws.Cells.Item(150, 4).NumberFormat
